# Tabraiz Shamsi.xlsx edit
# 1. Insert a new column at the very left (A), shifting the existing
#    teamName..result columns from A:L to B:M.
# 2. Populate the new column A with header "matchNo" and the match-number
#    value "36th" for the single data row.
# 3. Rename the sheet from the default "Sheet1" to "Tabraiz Shamsi".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; all data (A:L) shifts to (B:M).
$ws.Columns.Item(1).Insert()

# New column A content.
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "36th"

# Rename the worksheet tab to match the player name.
$ws.Name = "Tabraiz Shamsi"
